$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new log entries for rows 10-14: dates in column A (formatted like the
# existing A9 date entry) and a logged hour (1) in column H (Self-guided
# learning hours / Other).

$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = 45628
$ws.Range("H10").Value = 1

$ws.Range("A9").Copy($ws.Range("A11"))
$ws.Range("A11").Value = 45629
$ws.Range("H11").Value = 1

$ws.Range("A9").Copy($ws.Range("A12"))
$ws.Range("A12").Value = 45630
$ws.Range("H12").Value = 1

$ws.Range("A9").Copy($ws.Range("A13"))
$ws.Range("A13").Value = 45635
$ws.Range("H13").Value = 1

$ws.Range("A9").Copy($ws.Range("A14"))
$ws.Range("A14").Value = 45636
$ws.Range("H14").Value = 1

# Move the selection/active cell to A15, matching the saved cursor position.
$ws.Range("A15").Select()
